# Update the workbook from FHIR StructureDefinition version 5.0.0 to 6.0.0
# (Alvearie/alvearie-fhir-ig deploy), per commit:
#   "Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450..."
#
# Changes:
#  - Metadata sheet: Version 5.0.0 -> 6.0.0; Date bumped; Publisher set to
#    "Alvearie Team"; "Contact / No display for ContactDetail" row replaced
#    by "Jurisdiction / United States of America"; duplicate Contact row removed.
#  - Elements sheet: root Extension row's Short/Definition updated to the
#    new benefit-plan wording.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# Replace the "Contact" row with "Jurisdiction"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact / No display for ContactDetail" row that
# no longer exists - delete the entire row, shifting rows 12+ up by one.
$meta.Range("A11:B11").EntireRow.Delete()

# --- Elements sheet updates ---
# Root Extension element's Short (K2) and Definition (L2) text changed.
$elements.Range("K2").Value = "Claim Response Benefit Plan"
$elements.Range("L2").Value = "The benefit plan on the claim item"
